# Update crypto price/volume data (and swap HuobiToken/MXToken rows 36-37)
# to match the refreshed GitHub Actions data pull.
# Note: numeric-looking price strings (e.g. "1.001") are prefixed with a
# leading apostrophe so Excel keeps them as text instead of parsing them
# as numbers (matching the original inlineStr/text storage).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.081.58'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.633.75'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '''214.06'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = '''0.5215'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '''0.2600'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('D9').Value = '''0.06286'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').Value = '''20.56'
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('D11').Value = '''0.07615'
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').Value = '1.650.19'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').Value = '''4.415'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').Value = '1.858.36'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '''0.5490'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').Value = '0.0₅8045'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').Value = '''64.80'
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('D18').Value = '26.053.14'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '''4.673'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').Value = '''188.34'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').Value = '''10.15'
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('D23').Value = '''6.128'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('D24').Value = '''1.002'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').Value = '''145.60'
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('D27').Value = '''7.395'
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('D28').Value = '''15.83'
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').Value = '''1.390'
$ws.Range('E29').Value = '  +3.03%  '
$ws.Range('D30').Value = '''0.05860'
$ws.Range('E30').Value = '  -6.41%  '
$ws.Range('D31').Value = '''1.258'
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('D32').Value = '''3.431'
$ws.Range('E32').Value = '  -1.90%  '
$ws.Range('D33').Value = '''3.397'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D34').Value = '''1.637'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '''0.9823'
$ws.Range('E35').Value = '  -1.37%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '''2.765'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '''2.397'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').Value = '''0.5738'
$ws.Range('E38').Value = '  -5.07%  '
$ws.Range('D39').Value = '''0.01615'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').Value = '''0.8566'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('D41').Value = '''1.001'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').Value = '1.032.82'
$ws.Range('E42').Value = '  -6.57%  '
$ws.Range('D43').Value = '''5.652'
$ws.Range('E43').Value = '  -7.72%  '
$ws.Range('D44').Value = '''100.24'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '1.784.33'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '''55.34'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').Value = '''8.090'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').Value = '''0.9982'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('D50').Value = '''0.05166'
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('D51').Value = '''0.4222'
$ws.Range('E51').Value = '  -0.75%  '
